$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 9 (un_franzosa_ControlvsCD_Fp)
# to make room for the new "un_franzosa_ControlvsCD_ConvCD" entry.
$ws.Rows("9:9").Insert()

# Insert a new row before the current row 14 (un_franzosa_ControlvsUC_Fp, after the
# first insertion shifted it down from row 13) to make room for the new
# "un_franzosa_ControlvsUC_ConvUC" entry.
$ws.Rows("15:15").Insert()

# Rewrite rows 9-28 with the final data (existing rows shifted plus the two new rows).
$ws.Cells.Item(9, 1).Value = "un_franzosa_ControlvsCD_ConvCD"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.33
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.67
$ws.Cells.Item(9, 8).Value = 0.67
$ws.Cells.Item(10, 1).Value = "un_franzosa_ControlvsCD_Fp"
$ws.Cells.Item(10, 2).Value = 0.33
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0.67
$ws.Cells.Item(10, 7).Value = 0.67
$ws.Cells.Item(10, 8).Value = 0.67
$ws.Cells.Item(11, 1).Value = "un_franzosa_ControlvsDisease_Age"
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0.33
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.67
$ws.Cells.Item(11, 8).Value = 0.67
$ws.Cells.Item(12, 1).Value = "un_franzosa_ControlvsDisease_ConvDisease"
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0.67
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.33
$ws.Cells.Item(12, 8).Value = 0.33
$ws.Cells.Item(13, 1).Value = "un_franzosa_ControlvsDisease_Fp"
$ws.Cells.Item(13, 2).Value = 0.33
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.33
$ws.Cells.Item(13, 6).Value = 0.67
$ws.Cells.Item(13, 7).Value = 0.33
$ws.Cells.Item(13, 8).Value = 0.33
$ws.Cells.Item(14, 1).Value = "un_franzosa_ControlvsUC_Age"
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 1
$ws.Cells.Item(15, 1).Value = "un_franzosa_ControlvsUC_ConvUC"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0.33
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.67
$ws.Cells.Item(15, 8).Value = 0.67
$ws.Cells.Item(16, 1).Value = "un_franzosa_ControlvsUC_Fp"
$ws.Cells.Item(16, 2).Value = 0.33
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0.33
$ws.Cells.Item(16, 6).Value = 0.67
$ws.Cells.Item(16, 7).Value = 0.33
$ws.Cells.Item(16, 8).Value = 0.33
$ws.Cells.Item(17, 1).Value = "nf_yachida_age"
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0.33
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.67
$ws.Cells.Item(17, 8).Value = 0.67
$ws.Cells.Item(18, 1).Value = "nf_yachida_alcohol"
$ws.Cells.Item(18, 2).Value = 0
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0.33
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.67
$ws.Cells.Item(18, 8).Value = 0.67
$ws.Cells.Item(19, 1).Value = "nf_yachida_BrinkmanIndex"
$ws.Cells.Item(19, 2).Value = 0
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(20, 1).Value = "nf_yachida_gender"
$ws.Cells.Item(20, 2).Value = 0
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0.67
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.33
$ws.Cells.Item(20, 8).Value = 0.33
$ws.Cells.Item(21, 1).Value = "nf_yachida_healthyvscancer"
$ws.Cells.Item(21, 2).Value = 0
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0.33
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.67
$ws.Cells.Item(21, 8).Value = 0.67
$ws.Cells.Item(22, 1).Value = "nf_yachida_healthyvsstageIII_IV"
$ws.Cells.Item(22, 2).Value = 0
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0.67
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 0.33
$ws.Cells.Item(22, 8).Value = 0.33
$ws.Cells.Item(23, 1).Value = "nf_wang_age"
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0.67
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 0.33
$ws.Cells.Item(23, 8).Value = 0.33
$ws.Cells.Item(24, 1).Value = "nf_wang_bmi"
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 0.33
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 0.67
$ws.Cells.Item(24, 8).Value = 0.67
$ws.Cells.Item(25, 1).Value = "nf_wang_creatinine"
$ws.Cells.Item(25, 2).Value = 0.33
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0.67
$ws.Cells.Item(25, 6).Value = 0.67
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(26, 1).Value = "nf_wang_egfr"
$ws.Cells.Item(26, 2).Value = 0
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 0.67
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 0.33
$ws.Cells.Item(26, 8).Value = 0.33
$ws.Cells.Item(27, 1).Value = "nf_wang_studygroup"
$ws.Cells.Item(27, 2).Value = 0
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 0.67
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 0.33
$ws.Cells.Item(27, 8).Value = 0.33
$ws.Cells.Item(28, 1).Value = "nf_wang_urea"
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0.67
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 0.33
$ws.Cells.Item(28, 8).Value = 0.33
